$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix punctuation in "Razon social" (column E) entries where commas
#        (and, in one case, internal dots in an abbreviation) were mangled
#        during a prior scrape/import. ---

$nameFixes = @{
    "E65"  = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN";
    "E225" = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN";
    "E117" = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO";
    "E215" = "RICCOTTI. MARIANA EDITH";
    "E240" = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH";
}

foreach ($addr in $nameFixes.Keys) {
    $ws.Range($addr).Value = $nameFixes[$addr]
}

# --- 2) Re-format every "Importe" value (column H, rows 2-322) from
#        Argentine style (thousands separator "." and decimal comma ",")
#        into a plain decimal string using "." as the decimal separator,
#        while keeping the cell stored as text (not coerced into a real
#        floating point number, which would also lose trailing zeros and
#        introduce binary rounding noise). ---

$importeRange = $ws.Range("H2:H322")
$importeRange.NumberFormat = "@"

for ($r = 2; $r -le 322; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $old = $cell.Text
    if ($old -ne "") {
        $new = $old.Replace(".", "").Replace(",", ".")
        if ($new -ne $old) {
            $cell.Value = $new
        }
    }
}

$importeRange.ClearFormats()
